# release app version 5.0
# Add 6 new customer rows (STT 53..58) to the KHACH_HANG sheet, mirroring
# existing rows already present in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KHACH_HANG")

function Set-EmptyTextCell {
    param($Sheet, [int]$Row, [int]$Col)
    # A bare apostrophe forces a genuine empty *text* cell (quote-prefix
    # empty string) instead of Excel's normal "clear the cell" behaviour
    # that a plain Value = "" triggers.
    $Sheet.Cells.Item($Row, $Col).Value = "'"
}

function Set-GeneralTextCell {
    # For columns whose default column style has NO explicit number format
    # (General), Excel's COM layer auto-detects date-looking strings (e.g.
    # "08/03/1995") and silently converts them to date serials. Temporarily
    # forcing a Text format avoids that, and switching back to a
    # (lower-case) "general" format afterwards maps back onto the sheet's
    # original plain style instead of minting a new one.
    param($Sheet, [int]$Row, [int]$Col, [string]$Text)
    if ([string]::IsNullOrEmpty($Text)) {
        Set-EmptyTextCell $Sheet $Row $Col
    } else {
        $Sheet.Cells.Item($Row, $Col).NumberFormat = "@"
        $Sheet.Cells.Item($Row, $Col).Value = $Text
        $Sheet.Cells.Item($Row, $Col).NumberFormat = "general"
    }
}

function Set-PlainCell {
    # Columns that already carry an explicit Text number format (@) at the
    # column level never get auto-converted, so a direct Value assignment
    # round-trips exactly as typed.
    param($Sheet, [int]$Row, [int]$Col, [string]$Text)
    if ([string]::IsNullOrEmpty($Text)) {
        Set-EmptyTextCell $Sheet $Row $Col
    } else {
        $Sheet.Cells.Item($Row, $Col).Value = $Text
    }
}

$rows = @(
    @{ Row=54; STT=53; Name="PHAN HUY HOÀNG"; DOB="08/03/1995"; Phone="0935003445"; Email="huyhoang.phan@gmail.com"; IdNo="048095000369"; IdDate="13/04/2021"; IdPlace="Cục cảnh sát QLHC về TTXH"; Address="160/63 Trần Cao Vân, Phường Tam Thuận, Quận Thanh Khê, Thành phố Đà Nẵng"; Account="1199399699"; Bank="MB Bank"; ShipAddr="" },
    @{ Row=55; STT=54; Name="PHAN THANH ĐẠI"; DOB=""; Phone="0868134747"; Email="monstermax0007@gmail.com"; IdNo="066096019550"; IdDate="10/03/2023"; IdPlace="Cục cảnh sát QLHC về TTXH"; Address="Thôn 3 Ea Kao, Buôn Ma Thuột, Đắk Lắk"; Account="0231000668635"; Bank="Ngân hàng Vietcombank"; ShipAddr="" },
    @{ Row=56; STT=55; Name="TRẦN VĂN VƯƠNG"; DOB=""; Phone="0335640447"; Email="tranvuong46923@gmail.com"; IdNo="066201010148"; IdDate="09/09/2022"; IdPlace="Cục CSQLHC về TTXH"; Address="Buôn Jok, Ea H’đing, Cư M’gar, Đắk Lắk"; Account="0335640447"; Bank="VP Bank"; ShipAddr="" },
    @{ Row=57; STT=56; Name="PHAN HUY HOÀNG"; DOB="08/03/1995"; Phone="0935003445"; Email="huyhoang.phan@gmail.com"; IdNo="048095000369"; IdDate="13/04/2021"; IdPlace="Cục cảnh sát QLHC về TTXH"; Address="160/63 Trần Cao Vân, Phường Tam Thuận, Quận Thanh Khê, Thành phố Đà Nẵng"; Account="1199399699"; Bank="MB Bank"; ShipAddr="" },
    @{ Row=58; STT=57; Name="TRẦN THỊ HOÀNG LY"; DOB=""; Phone="0799345489"; Email="hoanglybds@gmail.com"; IdNo="048190006314"; IdDate="11/11/2021"; IdPlace="Cục cảnh sát QLHC về TTXH"; Address="Tổ 24, An Hải Bắc, Sơn Trà, Đà Nẵng"; Account=""; Bank=""; ShipAddr="" },
    @{ Row=59; STT=58; Name="TRẦN VĂN VƯƠNG"; DOB=""; Phone="0335640447"; Email="tranvuong46923@gmail.com"; IdNo="066201010148"; IdDate="09/09/2022"; IdPlace="Cục CSQLHC về TTXH"; Address="Buôn Jok, Ea H’đing, Cư M’gar, Đắk Lắk"; Account="0335640447"; Bank="VP Bank"; ShipAddr="" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.STT                        # A: STT (number)
    Set-PlainCell         $ws $row 2  $r.Name                     # B: Tên khách hàng (text col)
    Set-GeneralTextCell   $ws $row 3  $r.DOB                      # C: Ngày sinh (General col, date-like text)
    Set-PlainCell         $ws $row 4  $r.Phone                    # D: Số điện thoại (Text col)
    Set-PlainCell         $ws $row 5  $r.Email                    # E: Email (text col)
    Set-PlainCell         $ws $row 6  $r.IdNo                     # F: Số căn cước (Text col)
    Set-PlainCell         $ws $row 7  $r.IdDate                   # G: Ngày cấp (Text col, date-like but safe)
    Set-PlainCell         $ws $row 8  $r.IdPlace                  # H: Nơi cấp (text col)
    Set-PlainCell         $ws $row 9  $r.Address                  # I: Địa chỉ (text col)
    Set-PlainCell         $ws $row 10 $r.Account                  # J: Số tài khoản (Text col)
    Set-PlainCell         $ws $row 11 $r.Bank                     # K: Ngân hàng (text col)
    Set-PlainCell         $ws $row 12 $r.ShipAddr                 # L: Địa chỉ giao hàng (text col)
}
